$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.260.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.543.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.538.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  +3.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.18%  "

$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.37"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.122.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "609.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.546.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.315.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -17.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  -3.29%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("E29").Value = "  +3.54%  "

$ws.Range("E30").Value = "  -2.21%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.68%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "669.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.99%  "

$ws.Range("E33").Value = "  -3.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0474"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.373.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.85%  "

$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.34%  "

$ws.Range("E51").Value = "  -0.08%  "
